$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-text price strings (e.g. "30.741.66") that look
# numeric to Excel's automatic type detection. Force the range to Text
# format before writing so the values stay literal strings, then restore
# the original (default/"Normal") cell style once the writes are done.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.741.66"
$ws.Range("E2").Value = "  +2.67%  "
$ws.Range("D3").Value = "1.894.07"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "246.11"
$ws.Range("E5").Value = "  +2.04%  "
$ws.Range("D6").Value = "0.9997"
$ws.Range("D7").Value = "0.4932"
$ws.Range("E7").Value = "  -0.96%  "
$ws.Range("D8").Value = "0.2957"
$ws.Range("E8").Value = "  +1.24%  "
$ws.Range("D9").Value = "0.06813"
$ws.Range("E9").Value = "  +3.02%  "
$ws.Range("D10").Value = "17.42"
$ws.Range("E10").Value = "  +4.01%  "
$ws.Range("D11").Value = "1.891.07"
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("D12").Value = "92.69"
$ws.Range("E12").Value = "  +7.66%  "
$ws.Range("D13").Value = "0.07279"
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("D14").Value = "0.6837"
$ws.Range("E14").Value = "  +2.49%  "
$ws.Range("D15").Value = "5.090"
$ws.Range("E15").Value = "  +4.54%  "
$ws.Range("D16").Value = "30.715.47"
$ws.Range("E16").Value = "  +2.72%  "
$ws.Range("D17").Value = "0.000008008"
$ws.Range("E17").Value = "  +1.40%  "
$ws.Range("D18").Value = "13.32"
$ws.Range("D19").Value = "0.9993"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").Value = "2.139.97"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("D21").Value = "0.9942"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("E22").Value = "  +2.02%  "
$ws.Range("D23").Value = "193.65"
$ws.Range("E23").Value = "  +39.16%  "
$ws.Range("D24").Value = "6.102"
$ws.Range("E24").Value = "  +8.48%  "
$ws.Range("D25").Value = "9.336"
$ws.Range("E25").Value = "  +2.99%  "
$ws.Range("D26").Value = "155.02"
$ws.Range("E26").Value = "  +4.23%  "
$ws.Range("D27").Value = "19.64"
$ws.Range("E27").Value = "  +15.57%  "
$ws.Range("D28").Value = "1.927"
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("D29").Value = "1.391"
$ws.Range("E29").Value = "  +0.60%  "
$ws.Range("D30").Value = "4.345"
$ws.Range("E30").Value = "  +3.64%  "
$ws.Range("D31").Value = "0.09032"
$ws.Range("E31").Value = "  +2.76%  "
$ws.Range("D32").Value = "4.037"
$ws.Range("E32").Value = "  +2.11%  "
$ws.Range("E33").Value = "  +2.87%  "
$ws.Range("D34").Value = "0.7490"
$ws.Range("E34").Value = "  +5.71%  "
$ws.Range("D35").Value = "1.130"
$ws.Range("E35").Value = "  +2.12%  "
$ws.Range("D36").Value = "2.697"
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("E37").Value = "  +8.17%  "
$ws.Range("D38").Value = "2.675"
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("D39").Value = "2.167"
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("D40").Value = "0.9374"
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("D41").Value = "0.4454"
$ws.Range("E41").Value = "  +4.43%  "
$ws.Range("D42").Value = "105.88"
$ws.Range("E42").Value = "  +4.13%  "
$ws.Range("D43").Value = "5.838"
$ws.Range("E43").Value = "  +1.23%  "
$ws.Range("D44").Value = "0.9996"
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("D45").Value = "7.708"
$ws.Range("E45").Value = "  +3.83%  "
$ws.Range("D46").Value = "0.1342"
$ws.Range("E46").Value = "  +7.09%  "
$ws.Range("D47").Value = "0.05858"
$ws.Range("E47").Value = "  +3.57%  "
$ws.Range("D48").Value = "8.748"
$ws.Range("E48").Value = "  +6.34%  "
$ws.Range("D49").Value = "0.3974"
$ws.Range("D50").Value = "33.58"
$ws.Range("E50").Value = "  +4.02%  "
$ws.Range("D51").Value = "1.404"
$ws.Range("E51").Value = "  +5.40%  "

# Restore the default style on column D so no stray number-format is left
# applied to the cells (matches the source workbook, which used the
# default/unstyled cell format for these price cells).
$dRange.Style = "Normal"
